$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1630.7241
$ws.Range("I40").Value = 1517.7727
$ws.Range("K40").Value = 1517.7727
$ws.Range("M40").Value = -1342.7727
$ws.Range("H51").Value = 7324.9
$ws.Range("J51").Value = 2934.111
$ws.Range("L51").Value = 2934.111
$ws.Range("N51").Value = -3902.111
$ws.Range("H137").Value = 1489.6842
$ws.Range("I137").Value = 1344.24
$ws.Range("J137").Value = 1603.3125
$ws.Range("K137").Value = 4032.72
$ws.Range("L137").Value = 4809.9375
$ws.Range("M137").Value = -1482.72
$ws.Range("N137").Value = -9909.9375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 23067.662
$ws.Range("I32").Value = 4276.8486
$ws.Range("J32").Value = 111652.93
$ws.Range("K32").Value = 4276.8486
$ws.Range("L32").Value = 111652.93
$ws.Range("M32").Value = -3989.8486
$ws.Range("N32").Value = -112226.93
$ws.Range("H63").Value = 2706.0715
$ws.Range("I63").Value = 2370.5
$ws.Range("J63").Value = 3545
$ws.Range("K63").Value = 2370.5
$ws.Range("L63").Value = 3545
$ws.Range("M63").Value = -1684.5
$ws.Range("N63").Value = -4917
$ws.Range("H66").Value = 2706.0715
$ws.Range("I66").Value = 2370.5
$ws.Range("J66").Value = 3545
$ws.Range("K66").Value = 11852.5
$ws.Range("L66").Value = 17725
$ws.Range("M66").Value = -8420.5
$ws.Range("N66").Value = -24589
$ws.Range("H74").Value = 2148.5
$ws.Range("I74").Value = 1370.6296
$ws.Range("J74").Value = 3383.9412
$ws.Range("K74").Value = 1370.6296
$ws.Range("L74").Value = 3383.9412
$ws.Range("M74").Value = -496.6296
$ws.Range("N74").Value = -5131.9412
$ws.Range("H77").Value = 2148.5
$ws.Range("I77").Value = 1370.6296
$ws.Range("J77").Value = 3383.9412
$ws.Range("K77").Value = 6853.148
$ws.Range("L77").Value = 16919.706
$ws.Range("M77").Value = -2485.148
$ws.Range("N77").Value = -25655.706
$ws.Range("H132").Value = 2485.8594
$ws.Range("J132").Value = 1814.6
$ws.Range("L132").Value = 5443.799999999999
$ws.Range("N132").Value = -10503.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H27").Value = 29497.334
$ws.Range("H86").Value = 42486.594
$ws.Range("J86").Value = 1927.8889
$ws.Range("L86").Value = 1927.8889
$ws.Range("N86").Value = -4173.8889
$ws.Range("H89").Value = 42486.594
$ws.Range("J89").Value = 1927.8889
$ws.Range("L89").Value = 9639.4445
$ws.Range("N89").Value = -20871.4445

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 17858.484
$ws.Range("I31").Value = 56955.332
$ws.Range("J31").Value = 2559.7173
$ws.Range("K31").Value = 56955.332
$ws.Range("L31").Value = 2559.7173
$ws.Range("M31").Value = -56660.332
$ws.Range("N31").Value = -3149.7173
$ws.Range("H34").Value = 17858.484
$ws.Range("I34").Value = 56955.332
$ws.Range("J34").Value = 2559.7173
$ws.Range("K34").Value = 56955.332
$ws.Range("L34").Value = 2559.7173
$ws.Range("M34").Value = -56753.332
$ws.Range("N34").Value = -2963.7173
$ws.Range("H62").Value = 7938679.5
$ws.Range("J62").Value = 2533.3333
$ws.Range("L62").Value = 2533.3333
$ws.Range("N62").Value = -3781.3333
$ws.Range("H65").Value = 7938679.5
$ws.Range("J65").Value = 2533.3333
$ws.Range("L65").Value = 12666.6665
$ws.Range("N65").Value = -18906.6665
$ws.Range("H105").Value = 1061.5714
$ws.Range("I105").Value = 1021.8333
$ws.Range("J105").Value = 1300
$ws.Range("K105").Value = 1021.8333
$ws.Range("L105").Value = 1300
$ws.Range("M105").Value = 725.1667
$ws.Range("N105").Value = -4794
$ws.Range("H132").Value = 2261.9678
$ws.Range("I132").Value = 1996.36
$ws.Range("K132").Value = 5989.08
$ws.Range("M132").Value = -3459.08

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 6854.1665
$ws.Range("I3").Value = 3686.6667
$ws.Range("J3").Value = 10021.667
$ws.Range("K3").Value = 11060.0001
$ws.Range("L3").Value = 30065.001
$ws.Range("M3").Value = -10948.0001
$ws.Range("N3").Value = -30289.001
$ws.Range("H50").Value = 539.8461
$ws.Range("I50").Value = 347.8
$ws.Range("J50").Value = 659.875
$ws.Range("K50").Value = 1043.4
$ws.Range("L50").Value = 1979.625
$ws.Range("M50").Value = -562.4000000000001
$ws.Range("N50").Value = -2941.625
$ws.Range("H53").Value = 539.8461
$ws.Range("I53").Value = 347.8
$ws.Range("J53").Value = 659.875
$ws.Range("K53").Value = 1043.4
$ws.Range("L53").Value = 1979.625
$ws.Range("M53").Value = -562.4000000000001
$ws.Range("N53").Value = -2941.625
$ws.Range("H68").Value = 1966.7894
$ws.Range("I68").Value = 1404.6786
$ws.Range("J68").Value = 2294.6875
$ws.Range("K68").Value = 4214.0358
$ws.Range("L68").Value = 6884.0625
$ws.Range("M68").Value = -3403.0358
$ws.Range("N68").Value = -8506.0625
$ws.Range("H71").Value = 1966.7894
$ws.Range("I71").Value = 1404.6786
$ws.Range("J71").Value = 2294.6875
$ws.Range("K71").Value = 12642.1074
$ws.Range("L71").Value = 20652.1875
$ws.Range("M71").Value = -8586.107399999999
$ws.Range("N71").Value = -28764.1875
$ws.Range("H131").Value = 1516605.9
$ws.Range("I131").Value = 591.8
$ws.Range("J131").Value = 1710966.8
$ws.Range("K131").Value = 1775.4
$ws.Range("L131").Value = 5132900.4
$ws.Range("M131").Value = 3264.6
$ws.Range("N131").Value = -5142980.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H64").Value = 41565.8
$ws.Range("J64").Value = 41565.8
$ws.Range("L64").Value = 41565.8
$ws.Range("N64").Value = -42061.8
$ws.Range("H67").Value = 41565.8
$ws.Range("J67").Value = 41565.8
$ws.Range("L67").Value = 41565.8
$ws.Range("N67").Value = -43281.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 13167
$ws.Range("I132").Value = 29502
$ws.Range("J132").Value = 4999.5
$ws.Range("K132").Value = 88506
$ws.Range("L132").Value = 14998.5
$ws.Range("M132").Value = -85976
$ws.Range("N132").Value = -20058.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1395.9445
$ws.Range("I136").Value = 777.25
$ws.Range("J136").Value = 2633.3333
$ws.Range("K136").Value = 2331.75
$ws.Range("L136").Value = 7899.999899999999
$ws.Range("M136").Value = 218.25
$ws.Range("N136").Value = -12999.9999
